$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pasted-in lookup rows (original_id / replace_with pairs) for rows 2-10.
$rows = @(
    @(2,  "8ccf5f5e-e54f-43a8-af3f-cd2f367e417d", "f9f2e321-bc03-4c99-a096-763cfabac523"),
    @(3,  "b98d6333-5247-4e2b-b254-f55749f45ec9", "de868878-a260-46fe-8c9b-8997538a3c20"),
    @(4,  "0301e877-1470-40f2-adc6-43d199dec84c", "2fca9c78-69ab-4207-beac-e1d0a575e287"),
    @(5,  "cf8ff97e-6db5-4f45-8775-3bf1f29e411c", "2fca9c78-69ab-4207-beac-e1d0a575e287"),
    @(6,  "b6010b52-4672-4fe0-86dc-5cb488d6b383", "9c6a4d85-347d-44c3-95d0-ab3211db2fd4"),
    @(7,  "5598d097-021d-4c04-a466-01370403ebe9", "9bad9dca-5e02-4bd2-9fbc-1c268474d70c"),
    @(8,  "81a3cd60-939e-4306-962b-231012fdda70", "9bada080-fb4e-4b02-8f32-b1193fac3d5c"),
    @(9,  "c89d8952-34ba-4ebb-97dd-28708adf755d", "9bad9dcc-88c4-49e8-80bb-8ba040b1d89b"),
    @(10, "ea41f487-7fd5-4d00-bf20-9f72eb9056dd", "9bad9dcc-8e9e-406b-92db-ff709bba068f")
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
}

# Cursor ends up on A11 after filling the table, matching the saved selection.
$ws.Range("A11").Select()
